$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update the marking value (marks per correct answer) and recalculated totals
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 95
$ws.Range("E12").Value = "95/140"
